# Helper: write a value as TEXT (not auto-converted to a number) without
# leaving any stray NumberFormat / style behind. We do this by putting a
# text-producing formula into the cell, then converting it in place to a
# static value via Copy + PasteSpecial(xlPasteValues). This preserves the
# cell's original style (no new cellXf entries) while guaranteeing the
# stored type is a string (preserves leading zeros, avoids numeric coercion).
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new sheet "2022-Q3" right after "总计" and before "2022-Q1",
#    by duplicating the existing "2022-Q1" sheet (so it inherits the exact
#    same layout/styles) and then overwriting its data with the new figures.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsQ1.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Row 2 data
Set-TextValue $wsQ3.Range("B2") "501029"
Set-TextValue $wsQ3.Range("C2") "华宝标普中国A股红利机会指数（LOF）A"
Set-TextValue $wsQ3.Range("D2") "8.11"
Set-TextValue $wsQ3.Range("E2") "94.26"
Set-TextValue $wsQ3.Range("F2") "1.58"
Set-TextValue $wsQ3.Range("G2") "0.1281"
$wsQ3.Range("H2").Value = 7

# Row 3 data
Set-TextValue $wsQ3.Range("B3") "005125"
Set-TextValue $wsQ3.Range("C3") "华宝标普中国A股红利机会指数C"
Set-TextValue $wsQ3.Range("D3") "3.38"
Set-TextValue $wsQ3.Range("E3") "94.26"
Set-TextValue $wsQ3.Range("F3") "1.58"
Set-TextValue $wsQ3.Range("G3") "0.0534"
$wsQ3.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new top data row for
#    2022-Q3 and push the existing quarters down by one row.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Give the new bottom row (A5) the same style as the row above it (A4)
# before populating values, so it matches the existing "index column" look.
$wsTotal.Range("A4").Copy()
$wsTotal.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

# Shift old rows down by one (2020-Q4 -> row5, 2021-Q1 -> row4, 2022-Q1 -> row3)
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2020-Q4"
$wsTotal.Range("C5").Value = 5
$wsTotal.Range("D5").Value = 0.41

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q1"
$wsTotal.Range("C4").Value = 4
$wsTotal.Range("D4").Value = 2.64

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q1"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.22

# New top row: 2022-Q3
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.18
